# Crypto price/volume snapshot refresh (GitHub Actions scheduled update).
# Column D = Price, Column E = Volume(1h); both stored as literal text
# (not numbers) to match the feed formatting, e.g. "28.087.33" and "  -1.49%  ".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D): values that cannot be mistaken for a plain number ---
$ws.Range("D2").Value = "28.087.33"
$ws.Range("D3").Value = "1.898.39"
$ws.Range("D14").Value = "1.879.29"
$ws.Range("D23").Value = "28.146.60"
$ws.Range("D26").Value = "2.110.08"

# --- Price (column D): values that look numeric (e.g. "1.001", "41.80") -----
# Force Text number format first so Excel keeps the literal digits/trailing
# zeros instead of silently coercing the cell to a Double.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.34"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5034"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3900"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.80"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.414"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.303"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.46"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001110"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06637"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.228"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.46"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.312"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.89"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "158.06"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.56"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.085"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.622"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.585"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06602"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02400"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.305"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.224"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6477"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.979"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9993"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6115"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.41"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.304"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.688"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.004"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.29"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.191"

# --- Volume(1h) (column E): percentage strings, always safely stored as text ---
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("E8").Value = "  -1.42%  "
$ws.Range("E9").Value = "  -5.73%  "
$ws.Range("E10").Value = "  -2.92%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("E12").Value = "  -2.18%  "
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("E15").Value = "  -3.69%  "
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("E17").Value = "  -1.56%  "
$ws.Range("E18").Value = "  -2.92%  "
$ws.Range("E20").Value = "  -1.28%  "
$ws.Range("E22").Value = "  -1.15%  "
$ws.Range("E23").Value = "  -1.45%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  +1.40%  "
$ws.Range("E26").Value = "  -0.93%  "
$ws.Range("E27").Value = "  -6.70%  "
$ws.Range("E28").Value = "  -2.23%  "
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("E30").Value = "  -2.27%  "
$ws.Range("E31").Value = "  -1.98%  "
$ws.Range("E32").Value = "  -0.89%  "
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("E34").Value = "  -0.64%  "
$ws.Range("E35").Value = "  -3.25%  "
$ws.Range("E36").Value = "  -2.87%  "
$ws.Range("E37").Value = "  -1.97%  "
$ws.Range("E38").Value = "  -1.21%  "
$ws.Range("E39").Value = "  +9.24%  "
$ws.Range("E40").Value = "  -4.11%  "
$ws.Range("E41").Value = "  +0.30%  "
$ws.Range("E42").Value = "  -2.42%  "
$ws.Range("E43").Value = "  -1.97%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("E47").Value = "  +1.63%  "
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("E49").Value = "  -2.46%  "
$ws.Range("E51").Value = "  -1.71%  "
